$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ID")

# Remove the channel map entry on row 3 ("535" -> "flu").
# A3 keeps its style but loses its value; B3 is cleared entirely.
$ws.Range("A3:B3").ClearContents()

# Update the active selection to reflect the now-empty area below the table.
$ws.Activate()
$ws.Range("A6").Select()

$wb.Save()
